$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$ws.Range("A51").Value = "$ 27.391 CLP 30-10-20"
$ws.Range("A52").Value = "$ 27.391 CLP 30-10-20"
